$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("model")

# The 'type' column should directly reflect the recursively-expanded
# leaf type that matches the corresponding survey prompt, so the
# separate 'elementType' column is no longer needed.
$ws.Range("B4").Value = "geopoint"
$ws.Range("B5").Value = "select_one"

# Drop the now-redundant elementType column entirely.
$ws.Range("C1:C6").EntireColumn.Delete() | Out-Null

# The model sheet becomes the active tab/selection instead of survey.
$ws.Activate() | Out-Null
$ws.Range("C22").Select() | Out-Null
